$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = "Voto 'sim' de Jaques Wagner causa mal-estar e surpreende o STF"
$ws.Range("C2").Value = "https://g1.globo.com/politica/blog/natuza-nery/post/2023/11/23/voto-de-jaques-em-pec-racha-lideranca-do-governo-e-surpreende-ministros.ghtml"

# Update row 3
$ws.Range("A3").Value = "Hamas vai libertar amanhã 1º grupo de reféns com 13 mulheres e crianças"
$ws.Range("C3").Value = "https://g1.globo.com/mundo/noticia/2023/11/23/hamas-vai-libertar-13-civis-em-primeiro-grupo-diz-qatar.ghtml"

# Update row 4
$ws.Range("A4").Value = "Polícia encontra casa em SP com 33 pessoas treinadas para engolir drogas"
$ws.Range("C4").Value = "https://g1.globo.com/sp/sao-paulo/noticia/2023/11/23/homem-e-preso-suspeito-de-aliciar-e-treinar-mais-de-30-pessoas-para-engolir-drogas-e-leva-las-a-europa-grupo-foi-encontrado-com-passaportes-e-cocaina.ghtml"

# Update row 5
$ws.Range("A5").Value = "Dor terrível e até morte: o que pode acontecer com as 'mulas' do tráfico"
$ws.Range("C5").Value = "https://g1.globo.com/saude/noticia/2023/11/23/desmaios-dores-terriveis-convulsoes-e-morte-o-que-acontece-no-organismo-de-quem-transporta-drogas-no-estomago-no-reto-ou-na-vagina.ghtml"

# Update row 6
$ws.Range("A6").Value = "Celulares top de linha: g1 testa 4 modelos objetos de desejo"
$ws.Range("C6").Value = "https://g1.globo.com/guia/guia-de-compras/tecnologia/celulares/celulares-topo-de-linha-g1-testa-4-smartphones-que-sao-objetos-de-desejo.ghtml"

# Update row 7
$ws.Range("A7").Value = "EUA enviam documentos que comprovam recompra de relógio por advogado de Bolsonaro "
$ws.Range("C7").Value = "https://g1.globo.com/politica/blog/andreia-sadi/post/2023/11/23/pf-recebe-documentos-da-justica-americana-que-comprovam-recompra-de-relogio-de-bolsonaro-nos-eua.ghtml"

# Update row 8
$ws.Range("A8").Value = "Entenda proposta de federalizar a Cemig para pagar dívidas de MG"
$ws.Range("B8").Value = "Estado deve R$ 160 bilhões e negocia recuperação fiscal com o governo federal."
$ws.Range("C8").Value = "https://g1.globo.com/mg/minas-gerais/noticia/2023/11/23/entenda-a-divida-de-mg-com-a-uniao-o-que-e-o-regime-de-recuperacao-fiscal-e-a-alternativa-que-inclui-federalizacao-da-cemig.ghtml"

# Update row 9
$ws.Range("A9").Value = "Desabamento de prédio, gritos de socorro no mar e mais VÍDEOS do dia"
$ws.Range("C9").Value = "https://g1.globo.com/playlist/videos-para-assistir-agora.ghtml"

# Delete row 10 entirely (shift cells up)
$ws.Rows.Item(10).Delete()
